$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$ws.Range("B23").Value = "<50 Indstr. & Production$nl<25 Constr., Power engineering, Science, Education$nl<15 Transport, Trade, Services "
$ws.Range("B24").Value = "<100 Indstr. & Production$nl<50 Constr., Power engineering, Science, Education$nl<30 Transport, Trade, Services "
$ws.Range("B25").Value = ">=100 Indstr. & Production$nl>=50 Constr., Power engineering, Science, Education$nl>=30 Transport, Trade, Services "
